## Remove the last bullet item ("Conclusão") from the bulleted list that
## lives inside the table cell summarising the final report, merging the
## trailing "_GoBack" bookmark into the end of the preceding bullet
## ("Ficheiros Guardados"), exactly as shown by the source diff:
##
##   <w:r><w:t>Ficheiros Guardados</w:t></w:r>
##   <w:bookmarkStart w:id="0" w:name="_GoBack"/>
##   <w:bookmarkEnd w:id="0"/>
##   </w:p>
##
## (the whole "Conclusão" paragraph disappears).

$d = $word.ActiveDocument

# Locate the two bullet paragraphs we care about.
$ficheiros = $null
$conclusao = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Ficheiros Guardados*") {
        $ficheiros = $p
    }
    if ($t -like "*Conclus*") {
        $conclusao = $p
    }
}

# Position right before the paragraph mark that ends "Ficheiros Guardados".
$insertPos = $ficheiros.Range.End - 1

# A bookmark can't be collapsed exactly on a paragraph-mark boundary, so
# stage it around a throw-away placeholder character, then delete the
# placeholder -- the bookmark collapses cleanly to the boundary, leaving
# the "Ficheiros Guardados" run intact (not split in two).
$placeholder = $d.Range($insertPos, $insertPos)
$placeholder.InsertBefore("X")

$placeholderRange = $d.Range($insertPos, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange) | Out-Null

$placeholderRange2 = $d.Range($insertPos, $insertPos + 1)
$placeholderRange2.Delete()

# Finally, delete the entire "Conclusão" paragraph (text, paragraph mark
# and its original "_GoBack" bookmark all go away with it).
$conclusao.Range.Delete()
